# Add 2022-Q4 data: a new row on the "总计" summary sheet and a brand-new
# per-quarter detail sheet, inserted right after "总计" (shifting the other
# quarters' sheets/rows down).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "总计" (summary) sheet: insert the 2022-Q4 totals as the new row 2,
#    pushing the existing rows (2022-Q3 .. 2020-Q4) down by one row.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

for ($r = 8; $r -ge 2; $r--) {
    $dest = $r + 1
    $summary.Range("A$r").Copy($summary.Range("A$dest"))
    $summary.Range("B$dest").Value2 = $summary.Range("B$r").Value2
    $summary.Range("C$dest").Value2 = $summary.Range("C$r").Value2
    $summary.Range("D$dest").Value2 = $summary.Range("D$r").Value2
}

$summary.Range("A2").Value2 = 0
$summary.Range("B2").Value2 = "2022-Q4"
$summary.Range("C2").Value2 = 16
$summary.Range("D2").Value2 = 2.43

# ---------------------------------------------------------------------------
# 2. New "2022-Q4" detail sheet, positioned right after "总计".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($null, $summary)
$newSheet.Name = "2022-Q4"

# Header row.
$newSheet.Range("B1").Value2 = "基金代码"
$newSheet.Range("C1").Value2 = "基金名称"
$newSheet.Range("D1").Value2 = "基金规模"
$newSheet.Range("E1").Value2 = "股票总仓位"
$newSheet.Range("F1").Value2 = "仓位占比"
$newSheet.Range("G1").Value2 = "持有市值(亿元)"
$newSheet.Range("H1").Value2 = "仓位排名"

# Fund holdings rows (row index in column A is 0-based).
$data = @(
    @("010695", "华夏磐益一年定期开放混合",       "16.03", "98.69", "4.29", "0.6877", 7),
    @("009837", "华夏磐锐一年定期开放混合A",       "14.15", "75.21", "4.28", "0.6056", 3),
    @("161724", "招商中证煤炭等权指数（LOF）A",     "17.24", "93.84", "3.10", "0.5344", 4),
    @("217002", "招商安泰平衡混合",               "5.50",  "49.54", "3.50", "0.1925", 5),
    @("002317", "招商睿逸稳健配置混合",             "4.84",  "49.57", "3.42", "0.1655", 5),
    @("012964", "招商稳健平衡混合C",               "1.73",  "61.95", "3.45", "0.0597", 5),
    @("013596", "招商中证煤炭等权指数（LOF）C",     "1.56",  "93.84", "3.10", "0.0484", 4),
    @("012963", "招商稳健平衡混合A",               "1.02",  "61.95", "3.45", "0.0352", 5),
    @("013759", "招商精选平衡混合A",               "0.42",  "55.24", "4.91", "0.0206", 5),
    @("008736", "南方高股息主题股票A",             "0.79",  "91.27", "2.59", "0.0205", 9),
    @("014768", "景顺华城稳健6月持有混合C",         "1.61",  "22.61", "1.13", "0.0182", 4),
    @("009838", "华夏磐锐一年定期开放混合C",         "0.39",  "75.21", "4.28", "0.0167", 3),
    @("014767", "景顺华城稳健6月持有混合A",         "1.10",  "22.61", "1.13", "0.0124", 4),
    @("016347", "招商中证煤炭等权指数（LOF）E",     "0.20",  "93.84", "3.10", "0.0062", 4),
    @("013760", "招商精选平衡混合C",               "0.09",  "55.24", "4.91", "0.0044", 5),
    @("008737", "南方高股息主题股票C",             "0.07",  "91.27", "2.59", "0.0018", 9)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $item = $data[$i]

    $newSheet.Range("A$row").Value2 = $i
    # Columns B/D/E/F/G hold numeric-looking text (fund code, fund size,
    # position figures, ...) that must stay TEXT (matching the other quarter
    # sheets, and preserving fund codes' leading zeros) - the leading
    # apostrophe forces text storage instead of auto-converting to a number.
    $newSheet.Range("B$row").Value2 = "'" + $item[0]
    $newSheet.Range("C$row").Value2 = $item[1]
    $newSheet.Range("D$row").Value2 = "'" + $item[2]
    $newSheet.Range("E$row").Value2 = "'" + $item[3]
    $newSheet.Range("F$row").Value2 = "'" + $item[4]
    $newSheet.Range("G$row").Value2 = "'" + $item[5]
    $newSheet.Range("H$row").Value2 = $item[6]
}

# Pick up the same cell formatting (bold/bordered header row, centred index
# column) used by every other per-quarter sheet - "2022-Q3" (now shifted to
# position 3) is a same-shaped template.
$template = $wb.Worksheets.Item(3)
$template.Range("A1:H17").Copy()
$newSheet.Range("A1:H17").PasteSpecial(-4122)
